$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 with trade data, mirroring the structure of rows 2-12.
$ws.Range("A13").Value = 8965.26
$ws.Range("B13").Value = 9011.2199999999993
$ws.Range("C13").Value = 17.8
$ws.Range("D13").Value = 17.89
$ws.Range("E13").Value = $true
$ws.Range("F13").Value = 0.51
$ws.Range("G13").Value = 42620.766284722224
$ws.Range("H13").Value = $false

# Column G uses a date/time number format - reuse the existing style from
# the row above (instead of assigning NumberFormat text, which would create
# a brand-new custom numFmt) so the same style index is preserved.
$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$excel.CutCopyMode = 0
